# chore: update Sheets via scheduled runner
# Updates computed market-price / profit figures (columns H:N) for
# specific Leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# worksheets, reflecting refreshed pricing data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 810.44446
$ws.Range("I28").Value = 855.17645
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 855.17645
$ws.Range("L28").Value = 50
$ws.Range("M28").Value = -370.17645
$ws.Range("N28").Value = -1020

$ws.Range("H53").Value = 239.52632
$ws.Range("I53").Value = 127.583336
$ws.Range("J53").Value = 431.42856
$ws.Range("K53").Value = 127.583336
$ws.Range("L53").Value = 431.42856
$ws.Range("M53").Value = 509.416664
$ws.Range("N53").Value = -1705.42856

$ws.Range("H100").Value = 5627.727
$ws.Range("I100").Value = 7452.5
$ws.Range("J100").Value = 5222.222
$ws.Range("K100").Value = 7452.5
$ws.Range("L100").Value = 5222.222
$ws.Range("M100").Value = -6911.5
$ws.Range("N100").Value = -6304.222

$ws.Range("H129").Value = 1147.4314
$ws.Range("I129").Value = 409.25
$ws.Range("J129").Value = 1210.2554
$ws.Range("K129").Value = 1227.75
$ws.Range("L129").Value = 3630.7662
$ws.Range("M129").Value = 3772.25
$ws.Range("N129").Value = -13630.7662

$ws.Range("H132").Value = 1712.75
$ws.Range("I132").Value = 1374.7
$ws.Range("J132").Value = 3403
$ws.Range("K132").Value = 4124.1
$ws.Range("L132").Value = 10209
$ws.Range("M132").Value = -1594.1
$ws.Range("N132").Value = -15269


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 42479.043
$ws.Range("I2").Value = 800
$ws.Range("K2").Value = 800
$ws.Range("M2").Value = -687

$ws.Range("H5").Value = 1505.125
$ws.Range("I5").Value = 1648.7142
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 1648.7142
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -1536.7142
$ws.Range("N5").Value = -724

$ws.Range("H45").Value = 2152
$ws.Range("I45").Value = 1903
$ws.Range("J45").Value = 2650
$ws.Range("K45").Value = 1903
$ws.Range("L45").Value = 2650
$ws.Range("M45").Value = -1526
$ws.Range("N45").Value = -3404

$ws.Range("H116").Value = 42479.043
$ws.Range("I116").Value = 800
$ws.Range("K116").Value = 800
$ws.Range("M116").Value = 1494

$ws.Range("H118").Value = 62570.855
$ws.Range("J118").Value = 62570.855
$ws.Range("L118").Value = 62570.855
$ws.Range("N118").Value = -65884.85500000001

$ws.Range("H122").Value = 6584.579
$ws.Range("I122").Value = 6839.278
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 20517.834
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -18067.834
$ws.Range("N122").Value = -10900


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 42479.043
$ws.Range("I3").Value = 800
$ws.Range("K3").Value = 800
$ws.Range("M3").Value = -686

$ws.Range("H4").Value = 1505.125
$ws.Range("I4").Value = 1648.7142
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 1648.7142
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -1533.7142
$ws.Range("N4").Value = -730

$ws.Range("H64").Value = 393.75
$ws.Range("I64").Value = 210
$ws.Range("J64").Value = 700
$ws.Range("K64").Value = 210
$ws.Range("L64").Value = 700
$ws.Range("M64").Value = 15
$ws.Range("N64").Value = -1150

$ws.Range("H67").Value = 393.75
$ws.Range("I67").Value = 210
$ws.Range("J67").Value = 700
$ws.Range("K67").Value = 210
$ws.Range("L67").Value = 700
$ws.Range("M67").Value = 570
$ws.Range("N67").Value = -2260


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 258.33334
$ws.Range("I35").Value = 258.33334
$ws.Range("K35").Value = 258.33334
$ws.Range("M35").Value = 35.66665999999998

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H107").Value = 513.6923
$ws.Range("I107").Value = 507.25
$ws.Range("K107").Value = 507.25
$ws.Range("M107").Value = 1412.75

$ws.Range("H132").Value = 2917.383
$ws.Range("I132").Value = 2343.9062
$ws.Range("J132").Value = 4140.8
$ws.Range("K132").Value = 7031.7186
$ws.Range("L132").Value = 12422.4
$ws.Range("M132").Value = -4501.7186
$ws.Range("N132").Value = -17482.4


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 131.23077
$ws.Range("I14").Value = 131.23077
$ws.Range("K14").Value = 393.69231
$ws.Range("M14").Value = -220.69231

$ws.Range("H68").Value = 179566.05
$ws.Range("I68").Value = 455198.53
$ws.Range("J68").Value = 1215.6177
$ws.Range("K68").Value = 1365595.59
$ws.Range("L68").Value = 3646.8531
$ws.Range("M68").Value = -1364784.59
$ws.Range("N68").Value = -5268.8531

$ws.Range("H71").Value = 179566.05
$ws.Range("I71").Value = 455198.53
$ws.Range("J71").Value = 1215.6177
$ws.Range("K71").Value = 4096786.77
$ws.Range("L71").Value = 10940.5593
$ws.Range("M71").Value = -4092730.77
$ws.Range("N71").Value = -19052.5593

$ws.Range("H92").Value = 858
$ws.Range("J92").Value = 822.5
$ws.Range("L92").Value = 2467.5
$ws.Range("N92").Value = -4963.5


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3063.4546
$ws.Range("I132").Value = 2490.1177
$ws.Range("J132").Value = 3672.625
$ws.Range("K132").Value = 7470.353099999999
$ws.Range("L132").Value = 11017.875
$ws.Range("M132").Value = -4940.353099999999
$ws.Range("N132").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 650
$ws.Range("I17").Value = 650
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 650
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -480
$ws.Range("N17").ClearContents()

$ws.Range("H46").Value = 1018.7778
$ws.Range("I46").Value = 771.2857
$ws.Range("J46").Value = 1885
$ws.Range("K46").Value = 771.2857
$ws.Range("L46").Value = 1885
$ws.Range("M46").Value = -583.2857
$ws.Range("N46").Value = -2261


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 769250
$ws.Range("I20").Value = 769250
$ws.Range("K20").Value = 769250
$ws.Range("M20").Value = -769010

$ws.Range("H24").Value = 53006
$ws.Range("I24").Value = 27500
$ws.Range("J24").Value = 70010
$ws.Range("K24").Value = 27500
$ws.Range("L24").Value = 70010
$ws.Range("M24").Value = -27270
$ws.Range("N24").Value = -70470

$ws.Range("H107").Value = 408.63635
$ws.Range("I107").Value = 326.8
$ws.Range("K107").Value = 980.4000000000001
$ws.Range("M107").Value = 939.5999999999999

$ws.Range("H132").Value = 2232.7354
$ws.Range("I132").Value = 2048.2693
$ws.Range("J132").Value = 2832.25
$ws.Range("K132").Value = 6144.8079
$ws.Range("L132").Value = 8496.75
$ws.Range("M132").Value = -3614.8079
$ws.Range("N132").Value = -13556.75

$ws.Range("H136").Value = 2406.3142
$ws.Range("I136").Value = 2550.889
$ws.Range("J136").Value = 2253.2354
$ws.Range("K136").Value = 7652.667
$ws.Range("L136").Value = 6759.706200000001
$ws.Range("M136").Value = -5102.667
$ws.Range("N136").Value = -11859.7062

